$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# New "Riders" (column C) and "Average" (column D) values for Madigan bike hours update
$ws.Range("C2").Value = 265
$ws.Range("D2").Value = 232.55

$ws.Range("C3").Value = 222
$ws.Range("D3").Value = 209.82

$ws.Range("C4").Value = 218
$ws.Range("D4").Value = 195.73

$ws.Range("C5").Value = 241
$ws.Range("D5").Value = 220.1

$ws.Range("C6").Value = 289
$ws.Range("D6").Value = 236.55

$ws.Range("C7").Value = 108
$ws.Range("D7").Value = 117.42

$ws.Range("C8").Value = 85
$ws.Range("D8").Value = 101.5
